# Bastiglia.xlsx report update: a new daily observation (2021-02-08, serial 44235,
# nuovi pos. = 2) is inserted into the time series between the existing
# 2021-02-07 (row 92) and 2021-02-09 rows, shifting every following row down by
# one. Two brand-new trailing rows are also appended for 2021-03-01 and
# 2021-03-02. The rolling 7-day sum (column C) and the per-100k-inhabitants
# figure (column D) are recomputed for every affected row, matching the
# target workbook exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 93 -- this pushes the old rows 93..113 down
# to 94..114 automatically, carrying along their existing values/styles.
$ws.Rows.Item(93).Insert()

# The freshly inserted row's A cell doesn't carry the date style (bold,
# centered, bordered, custom date format) used by the rest of column A, so
# clone it from a neighboring, correctly-styled cell.
$ws.Cells.Item(92,1).Copy()
$ws.Cells.Item(93,1).PasteSpecial(-4122)

# Likewise, row 115 is brand new territory (beyond the old used range), so
# clone the date style onto its A cell too.
$ws.Cells.Item(114,1).Copy()
$ws.Cells.Item(115,1).PasteSpecial(-4122)

# Now (re)write the literal values for every row from 90 through the new
# last row, 115, matching the target workbook.
$ws.Cells.Item(90,1).Value = 44232
$ws.Cells.Item(90,2).Value = 1
$ws.Cells.Item(90,3).Value = 6
$ws.Cells.Item(90,4).Value = 142.9592566118656

$ws.Cells.Item(91,1).Value = 44233
$ws.Cells.Item(91,2).Value = 1
$ws.Cells.Item(91,3).Value = 5
$ws.Cells.Item(91,4).Value = 119.1327138432213

$ws.Cells.Item(92,1).Value = 44234
$ws.Cells.Item(92,2).Value = 0
$ws.Cells.Item(92,3).Value = 5
$ws.Cells.Item(92,4).Value = 119.1327138432213

$ws.Cells.Item(93,1).Value = 44235
$ws.Cells.Item(93,2).Value = 2
$ws.Cells.Item(93,3).Value = 4
$ws.Cells.Item(93,4).Value = 95.30617107457708

$ws.Cells.Item(94,1).Value = 44236
$ws.Cells.Item(94,2).Value = 0
$ws.Cells.Item(94,3).Value = 3
$ws.Cells.Item(94,4).Value = 71.47962830593281

$ws.Cells.Item(95,1).Value = 44237
$ws.Cells.Item(95,2).Value = 0
$ws.Cells.Item(95,3).Value = 2
$ws.Cells.Item(95,4).Value = 47.65308553728854

$ws.Cells.Item(96,1).Value = 44238
$ws.Cells.Item(96,2).Value = 0
$ws.Cells.Item(96,3).Value = 2
$ws.Cells.Item(96,4).Value = 47.65308553728854

$ws.Cells.Item(97,1).Value = 44239
$ws.Cells.Item(97,2).Value = 0
$ws.Cells.Item(97,3).Value = 2
$ws.Cells.Item(97,4).Value = 47.65308553728854

$ws.Cells.Item(98,1).Value = 44240
$ws.Cells.Item(98,2).Value = 0
$ws.Cells.Item(98,3).Value = 3
$ws.Cells.Item(98,4).Value = 71.47962830593281

$ws.Cells.Item(99,1).Value = 44241
$ws.Cells.Item(99,2).Value = 0
$ws.Cells.Item(99,3).Value = 3
$ws.Cells.Item(99,4).Value = 71.47962830593281

$ws.Cells.Item(100,1).Value = 44242
$ws.Cells.Item(100,2).Value = 2
$ws.Cells.Item(100,3).Value = 3
$ws.Cells.Item(100,4).Value = 71.47962830593281

$ws.Cells.Item(101,1).Value = 44243
$ws.Cells.Item(101,2).Value = 1
$ws.Cells.Item(101,3).Value = 5
$ws.Cells.Item(101,4).Value = 119.1327138432213

$ws.Cells.Item(102,1).Value = 44244
$ws.Cells.Item(102,2).Value = 0
$ws.Cells.Item(102,3).Value = 7
$ws.Cells.Item(102,4).Value = 166.7857993805099

$ws.Cells.Item(103,1).Value = 44245
$ws.Cells.Item(103,2).Value = 0
$ws.Cells.Item(103,3).Value = 11
$ws.Cells.Item(103,4).Value = 262.091970455087

$ws.Cells.Item(104,1).Value = 44246
$ws.Cells.Item(104,2).Value = 2
$ws.Cells.Item(104,3).Value = 10
$ws.Cells.Item(104,4).Value = 238.2654276864427

$ws.Cells.Item(105,1).Value = 44247
$ws.Cells.Item(105,2).Value = 2
$ws.Cells.Item(105,3).Value = 10
$ws.Cells.Item(105,4).Value = 238.2654276864427

$ws.Cells.Item(106,1).Value = 44248
$ws.Cells.Item(106,2).Value = 4
$ws.Cells.Item(106,3).Value = 10
$ws.Cells.Item(106,4).Value = 238.2654276864427

$ws.Cells.Item(107,1).Value = 44249
$ws.Cells.Item(107,2).Value = 1
$ws.Cells.Item(107,3).Value = 11
$ws.Cells.Item(107,4).Value = 262.091970455087

$ws.Cells.Item(108,1).Value = 44250
$ws.Cells.Item(108,2).Value = 1
$ws.Cells.Item(108,3).Value = 12
$ws.Cells.Item(108,4).Value = 285.9185132237312

$ws.Cells.Item(109,1).Value = 44251
$ws.Cells.Item(109,2).Value = 0
$ws.Cells.Item(109,3).Value = 10
$ws.Cells.Item(109,4).Value = 238.2654276864427

$ws.Cells.Item(110,1).Value = 44252
$ws.Cells.Item(110,2).Value = 1
$ws.Cells.Item(110,3).Value = 8
$ws.Cells.Item(110,4).Value = 190.6123421491542

$ws.Cells.Item(111,1).Value = 44253
$ws.Cells.Item(111,2).Value = 3
$ws.Cells.Item(111,3).Value = 12
$ws.Cells.Item(111,4).Value = 285.9185132237312

$ws.Cells.Item(112,1).Value = 44254
$ws.Cells.Item(112,2).Value = 0
$ws.Cells.Item(112,3).Value = 12
$ws.Cells.Item(112,4).Value = 285.9185132237312

$ws.Cells.Item(113,1).Value = 44255
$ws.Cells.Item(113,2).Value = 2
$ws.Cells.Item(113,3).ClearContents()
$ws.Cells.Item(113,4).ClearContents()

$ws.Cells.Item(114,1).Value = 44256
$ws.Cells.Item(114,2).Value = 5
$ws.Cells.Item(114,3).ClearContents()
$ws.Cells.Item(114,4).ClearContents()

$ws.Cells.Item(115,1).Value = 44257
$ws.Cells.Item(115,2).Value = 1
$ws.Cells.Item(115,3).ClearContents()
$ws.Cells.Item(115,4).ClearContents()
